$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several of the new "Price" values are plain decimal numbers (e.g. "0.634",
# "231.56", "1.00"). If assigned as a plain .Value, Excel auto-converts them
# to numeric cells, which loses the exact original text (e.g. trailing
# zeros, or scientific notation for tiny numbers). Force these specific
# cells to Text format first so the values round-trip as exact strings,
# matching the source data (the other Price cells, like "44.136.25", are
# not valid numbers so they remain text automatically).
# Note: Union() here only reliably combines up to 3 ranges per call, and
# setting a property directly on a multi-area Union result only affects its
# first area - so we build the combined range incrementally and then apply
# NumberFormat per-area via the Areas collection (this still only adds a
# single extra style to the workbook, reused by every cell).
$u = $ws.Range("D5")
$u = $excel.Union($u, $ws.Range("D6"), $ws.Range("D7"))
$u = $excel.Union($u, $ws.Range("D9"), $ws.Range("D10"))
$u = $excel.Union($u, $ws.Range("D11"), $ws.Range("D12"))
$u = $excel.Union($u, $ws.Range("D13"), $ws.Range("D15"))
$u = $excel.Union($u, $ws.Range("D16"), $ws.Range("D17"))
$u = $excel.Union($u, $ws.Range("D20"), $ws.Range("D21"))
$u = $excel.Union($u, $ws.Range("D22"), $ws.Range("D23"))
$u = $excel.Union($u, $ws.Range("D25"), $ws.Range("D27"))
$u = $excel.Union($u, $ws.Range("D28"), $ws.Range("D29"))
$u = $excel.Union($u, $ws.Range("D30"), $ws.Range("D31"))
$u = $excel.Union($u, $ws.Range("D32"), $ws.Range("D34"))
$u = $excel.Union($u, $ws.Range("D35"), $ws.Range("D36"))
$u = $excel.Union($u, $ws.Range("D37"), $ws.Range("D38"))
$u = $excel.Union($u, $ws.Range("D39"), $ws.Range("D41"))
$u = $excel.Union($u, $ws.Range("D42"), $ws.Range("D43"))
$u = $excel.Union($u, $ws.Range("D44"), $ws.Range("D45"))
$u = $excel.Union($u, $ws.Range("D46"), $ws.Range("D47"))
$u = $excel.Union($u, $ws.Range("D49"), $ws.Range("D50"))
$u = $excel.Union($u, $ws.Range("D51"), $ws.Range("D51"))
foreach ($area in $u.Areas) {
  $area.NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.982.05"
$ws.Range("E2").Value = "  +5.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.270.67"
$ws.Range("E3").Value = "  +3.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.42%  "

# Row 5 - BNB
$ws.Range("D5").Value = "231.56"
$ws.Range("E5").Value = "  +0.86%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +2.90%  "

# Row 7 - Solana
$ws.Range("D7").Value = "63.29"
$ws.Range("E7").Value = "  +4.87%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.429"
$ws.Range("E9").Value = "  +7.20%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +14.97%  "

# Row 11 - OKB
$ws.Range("D11").Value = "56.38"
$ws.Range("E11").Value = "  -1.02%  "

# Row 12 - now TRON (was Avalanche)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.106"
$ws.Range("E12").Value = "  +2.89%  "

# Row 13 - now Avalanche (was TRON)
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "25.82"
$ws.Range("E13").Value = "  +16.79%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.607.55"
$ws.Range("E14").Value = "  +3.03%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "15.69"
$ws.Range("E15").Value = "  +2.35%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "5.93"
$ws.Range("E16").Value = "  +6.59%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.823"
$ws.Range("E17").Value = "  +3.82%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.278.53"
$ws.Range("E18").Value = "  +3.39%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "43.967.27"
$ws.Range("E19").Value = "  +5.52%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  +15.12%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "73.62"
$ws.Range("E21").Value = "  +2.29%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "256.15"
$ws.Range("E23").Value = "  +6.19%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.01%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +4.98%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -6.11%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +4.65%  "

# Row 28 - Monero
$ws.Range("D28").Value = "171.34"
$ws.Range("E28").Value = "  +1.95%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "20.92"
$ws.Range("E29").Value = "  +6.29%  "

# Row 30 - Kaspa
$ws.Range("D30").Value = "0.138"
$ws.Range("E30").Value = "  -1.25%  "

# Row 31 - WEMIXToken
$ws.Range("D31").Value = "2.86"
$ws.Range("E31").Value = "  +10.31%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  -3.46%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +2.72%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0680"
$ws.Range("E34").Value = "  +5.86%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "4.71"
$ws.Range("E35").Value = "  +2.85%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "4.97"
$ws.Range("E36").Value = "  -0.33%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +6.63%  "

# Row 38 - THORChain
$ws.Range("D38").Value = "6.71"
$ws.Range("E38").Value = "  +7.30%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").Value = "2.33"
$ws.Range("E39").Value = "  +0.39%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +4.19%  "

# Row 41 - BinanceUSD
$ws.Range("D41").Value = "1.00"

# Row 42 - FraxShare
$ws.Range("D42").Value = "8.46"
$ws.Range("E42").Value = "  -2.66%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "17.46"
$ws.Range("E43").Value = "  +8.53%  "

# Row 44 - Cronos
$ws.Range("D44").Value = "0.0964"
$ws.Range("E44").Value = "  +1.18%  "

# Row 45 - FTXToken
$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  +1.29%  "

# Row 46 - Aave
$ws.Range("D46").Value = "97.54"
$ws.Range("E46").Value = "  +1.27%  "

# Row 47 - TrustWalletToken
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  -0.33%  "

# Row 48 - now Maker (was TerraClassic)
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.464.93"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49 - now NEARProtocol (was Maker)
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  +5.04%  "

# Row 50 - now TerraClassic (was NEARProtocol)
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").Value = "0.000205"
$ws.Range("E50").Value = "  -14.60%  "

# Row 51 - ARBITRUM
$ws.Range("D51").Value = "1.07"
$ws.Range("E51").Value = "  +0.82%  "
